$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Rosetta Sprint Planning")

# The only real data edit: "Completed Velocity" for 2016-06-06 (E10) drops from 6 to 0.
# All the Cumulative/Remaining columns (F, G, H, I) recompute automatically via their
# existing formulas once this single input changes.
$ws.Range("E10").Value = 0

$ws.Select()
$ws.Range("C8").Select()
